$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(32, 346),
    @(33, 301),
    @(34, 468),
    @(35, 711),
    @(36, 495),
    @(37, 509)
)

$row = 33
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
